$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.814.42'
$ws.Range('E2').Value = '  +0.07%  '

$ws.Range('D3').Value = '2.783.75'
$ws.Range('E3').Value = '  -1.21%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '357.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.17'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.21%  '

$ws.Range('E7').Value = '  -1.19%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.589'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.23%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.97'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.65%  '

$ws.Range('E11').Value = '  +2.13%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0846'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.53'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.71%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.61'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.84%  '

$ws.Range('D15').Value = '3.226.03'
$ws.Range('E15').Value = '  -0.99%  '

$ws.Range('D16').Value = '2.770.63'
$ws.Range('E16').Value = '  -2.34%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.940'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.42%  '

$ws.Range('D18').Value = '51.735.56'
$ws.Range('E18').Value = '  +0.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.45'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.69%  '

$ws.Range('E20').Value = '  -2.14%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.27'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.71%  '

$ws.Range('D22').Value = '0.0₃0973'
$ws.Range('E22').Value = '  -1.90%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.35'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.49%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.05'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.80%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.75'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.50%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.43'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.84%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.10%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.164'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +17.53%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.24'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.48%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.22'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.28%  '

$ws.Range('E31').Value = '  +6.73%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '52.03'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.18%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.92'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.67%  '

$ws.Range('E34').Value = '  -7.60%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0843'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.12%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.13'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -5.27%  '

$ws.Range('E37').Value = '  +0.02%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.74'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.33%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.15'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.28%  '

$ws.Range('E40').Value = '  -3.47%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.57'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.01%  '

$ws.Range('E42').Value = '  -1.93%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.22'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.11%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '119.81'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.46%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.71'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -6.55%  '

$ws.Range('D46').Value = '2.081.76'
$ws.Range('E46').Value = '  -0.63%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.28'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.64%  '

$ws.Range('E48').Value = '  +1.21%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.74'
$ws.Range('D49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.939'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.83%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.192'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.94%  '
